# Append a new log row (row 42) to each of the 4 worksheets, duplicating the
# last existing row (row 41) but with an updated timestamp in column A.
#
# Sheet mapping:
#   1 -> DE_LFT_#1
#   2 -> DE_LFT_#2
#   3 -> DE_PLT_#1
#   4 -> DE_PLT_#2

$wb = $excel.ActiveWorkbook

$newTime = 45828.43628472222
$dateFormat = "YYYY-MM-DD HH:MM:SS"

# Per-sheet column values for the new row 42 (copied from row 41, minus the
# timestamp in column A which advances to the new day).
$rowData = @{
    1 = @{
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x6C"
        E = "0x14"
        F = 380
        G = 759863127514710945038336.0
        H = 364
        I = 14
    }
    2 = @{
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x6C"
        E = "0xe"
        F = 380
        G = 568432987514711010443264.0
        H = 364
        I = 14
    }
    3 = @{
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x7F"
        E = "0x7"
        F = 130
        G = 568631262647113970876416.0
        H = 127
        I = 7
    }
    4 = @{
        B = "0x00,0x82"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x7F"
        E = "0x3"
        F = 130
        G = 985046333984776009023488.0
        H = 127
        I = 3
    }
}

foreach ($sheetIndex in 1..4) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    $data = $rowData[$sheetIndex]
    $row = 42

    $ws.Cells.Item($row, 1).Value = $newTime
    $ws.Cells.Item($row, 1).NumberFormat = $dateFormat

    $ws.Cells.Item($row, 2).Value = $data.B
    $ws.Cells.Item($row, 3).Value = $data.C
    $ws.Cells.Item($row, 4).Value = $data.D
    $ws.Cells.Item($row, 5).Value = $data.E
    $ws.Cells.Item($row, 6).Value = $data.F
    $ws.Cells.Item($row, 7).Value = $data.G
    $ws.Cells.Item($row, 8).Value = $data.H
    $ws.Cells.Item($row, 9).Value = $data.I
}
